$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulation results for the 380 kV case (pl_mw.xlsx): columns B-F and K-M
# for rows 2-25 are refreshed with the newly computed values.
$ws.Range("B2").Value = 1.092944669800971
$ws.Range("C2").Value = 0.0402866848161807
$ws.Range("D2").Value = 0.3469786856266808
$ws.Range("E2").Value = 0.08423663855832686
$ws.Range("F2").Value = 4.953339856236369
$ws.Range("K2").Value = 0.6178926860127945
$ws.Range("L2").Value = 0.1990288344996003
$ws.Range("M2").Value = 0.2564379338706004
$ws.Range("B3").Value = 1.079093574167189
$ws.Range("C3").Value = 0.03812766800084333
$ws.Range("D3").Value = 0.3343197433335092
$ws.Range("E3").Value = 0.08360673564671295
$ws.Range("F3").Value = 4.74137407832751
$ws.Range("K3").Value = 0.6045266538105096
$ws.Range("L3").Value = 0.196850744342612
$ws.Range("M3").Value = 0.2534131723321948
$ws.Range("B4").Value = 1.071552234929982
$ws.Range("C4").Value = 0.03677288371079257
$ws.Range("D4").Value = 0.32645662235484
$ws.Range("E4").Value = 0.08323112675766886
$ws.Range("F4").Value = 4.611407365189393
$ws.Range("K4").Value = 0.5970327252603767
$ws.Range("L4").Value = 0.1956125514581331
$ws.Range("M4").Value = 0.2517483755631105
$ws.Range("B5").Value = 1.068720938830609
$ws.Range("C5").Value = 0.03621330098894049
$ws.Range("D5").Value = 0.3232291781499015
$ws.Range("E5").Value = 0.08308086286111305
$ws.Range("F5").Value = 4.558486500939694
$ws.Range("K5").Value = 0.5941575649632682
$ws.Range("L5").Value = 0.1951329034367504
$ws.Range("M5").Value = 0.251118292323973
$ws.Range("B6").Value = 1.068265401658692
$ws.Range("C6").Value = 0.03611992494107596
$ws.Range("D6").Value = 0.3226918529031906
$ws.Range("E6").Value = 0.08305608061396974
$ws.Range("F6").Value = 4.549701432159594
$ws.Range("K6").Value = 0.5936909216756447
$ws.Range("L6").Value = 0.1950547634399413
$ws.Range("M6").Value = 0.251016585379805
$ws.Range("B7").Value = 1.071513072278407
$ws.Range("C7").Value = 0.03676536756045223
$ws.Range("D7").Value = 0.3264131902131879
$ws.Range("E7").Value = 0.08322908891594061
$ws.Range("F7").Value = 4.610693492545323
$ws.Range("K7").Value = 0.5969932271548117
$ws.Range("L7").Value = 0.1956059818456026
$ws.Range("M7").Value = 0.2517396823996449
$ws.Range("B8").Value = 1.087968648099263
$ws.Range("C8").Value = 0.0395482009699748
$ws.Range("D8").Value = 0.3426324133937584
$ws.Range("E8").Value = 0.08401712970679576
$ws.Range("F8").Value = 4.880214031325607
$ws.Range("K8").Value = 0.6131357622472962
$ws.Range("L8").Value = 0.1982572402477984
$ws.Range("M8").Value = 0.2553550265628353
$ws.Range("B9").Value = 1.127902360285674
$ws.Range("C9").Value = 0.04478189957499268
$ws.Range("D9").Value = 0.3737398177966611
$ws.Range("E9").Value = 0.08565131579503671
$ws.Range("F9").Value = 5.41038330910456
$ws.Range("K9").Value = 0.6504773667712129
$ws.Range("L9").Value = 0.2042441855646331
$ws.Range("M9").Value = 0.2639746959276508
$ws.Range("B10").Value = 1.161949957525053
$ws.Range("C10").Value = 0.04850182590726604
$ws.Range("D10").Value = 0.3961972016212201
$ws.Range("E10").Value = 0.08690674287358391
$ws.Range("F10").Value = 5.80121779326754
$ws.Range("K10").Value = 0.6814253061023408
$ws.Range("L10").Value = 0.2091253294511972
$ws.Range("M10").Value = 0.2712461558002559
$ws.Range("B11").Value = 1.178469861633317
$ws.Range("C11").Value = 0.05016929332256126
$ws.Range("D11").Value = 0.4063336508876887
$ws.Range("E11").Value = 0.08748991550168483
$ws.Range("F11").Value = 5.979376596327995
$ws.Range("K11").Value = 0.6962777568209333
$ws.Range("L11").Value = 0.2114512260868509
$ws.Range("M11").Value = 0.2747593367263761
$ws.Range("B12").Value = 1.184874437751745
$ws.Range("C12").Value = 0.05079735209854874
$ws.Range("D12").Value = 0.4101610952700412
$ws.Range("E12").Value = 0.08771249296683337
$ws.Range("F12").Value = 6.046898424744995
$ws.Range("K12").Value = 0.702014095431025
$ws.Range("L12").Value = 0.2123471769291285
$ws.Range("M12").Value = 0.27611931096353
$ws.Range("B13").Value = 1.183488469195026
$ws.Range("C13").Value = 0.05066223556759297
$ws.Range("D13").Value = 0.4093372694445065
$ws.Range("E13").Value = 0.08766447931481025
$ws.Range("F13").Value = 6.032353788764965
$ws.Range("K13").Value = 0.7007736776196225
$ws.Range("L13").Value = 0.2121535421919418
$ws.Range("M13").Value = 0.2758250980274966
$ws.Range("B14").Value = 1.178993783802525
$ws.Range("C14").Value = 0.0502210306777755
$ws.Range("D14").Value = 0.4066487552218518
$ws.Range("E14").Value = 0.08750819211957861
$ws.Range("F14").Value = 5.98493049151233
$ws.Range("K14").Value = 0.6967474392267832
$ws.Range("L14").Value = 0.2115246320540081
$ws.Range("M14").Value = 0.2748706287570002
$ws.Range("B15").Value = 1.176260059253792
$ws.Range("C15").Value = 0.04995034589655489
$ws.Range("D15").Value = 0.4050005416535214
$ws.Range("E15").Value = 0.08741268884481102
$ws.Range("F15").Value = 5.955889898372504
$ws.Range("K15").Value = 0.6942958644280566
$ws.Range("L15").Value = 0.2111413845585872
$ws.Range("M15").Value = 0.2742898471099053
$ws.Range("B16").Value = 1.160891143136894
$ws.Range("C16").Value = 0.04839237069239033
$ws.Range("D16").Value = 0.3955331969745259
$ws.Range("E16").Value = 0.08686887466292248
$ws.Range("F16").Value = 5.789582431684863
$ws.Range("K16").Value = 0.6804703016709368
$ws.Range("L16").Value = 0.2089754502136429
$ws.Range("M16").Value = 0.2710207001351534
$ws.Range("B17").Value = 1.15172736333048
$ws.Range("C17").Value = 0.04743039547035721
$ws.Range("D17").Value = 0.3897052218728732
$ws.Range("E17").Value = 0.08653835899190376
$ws.Range("F17").Value = 5.687655025184938
$ws.Range("K17").Value = 0.6721875383947236
$ws.Range("L17").Value = 0.2076737405945295
$ws.Range("M17").Value = 0.2690678362850321
$ws.Range("B18").Value = 1.146553655252035
$ws.Range("C18").Value = 0.04687475252167417
$ws.Range("D18").Value = 0.3863456135825771
$ws.Range("E18").Value = 0.08634939103524708
$ws.Range("F18").Value = 5.629063288590061
$ws.Range("K18").Value = 0.6674963274999186
$ws.Range("L18").Value = 0.2069349547281831
$ws.Range("M18").Value = 0.2679639249529586
$ws.Range("B19").Value = 1.144818580395906
$ws.Range("C19").Value = 0.04668621386225169
$ws.Range("D19").Value = 0.3852068049269093
$ws.Range("E19").Value = 0.0862856046287046
$ws.Range("F19").Value = 5.609230879382153
$ws.Range("K19").Value = 0.6659204483466965
$ws.Range("L19").Value = 0.2066865181756441
$ws.Range("M19").Value = 0.2675934761639311
$ws.Range("B20").Value = 1.152692814740817
$ws.Range("C20").Value = 0.04753304025356897
$ws.Range("D20").Value = 0.3903263940486568
$ws.Range("E20").Value = 0.08657342533820866
$ws.Range("F20").Value = 5.698501799924259
$ws.Range("K20").Value = 0.6730617113647099
$ws.Range("L20").Value = 0.2078112826488052
$ws.Range("M20").Value = 0.2692737214271261
$ws.Range("B21").Value = 1.180309937950057
$ws.Range("C21").Value = 0.05035071337431418
$ws.Range("D21").Value = 0.4074387323855149
$ws.Range("E21").Value = 0.08755405010843731
$ws.Range("F21").Value = 5.998858283204527
$ws.Range("K21").Value = 0.6979269964947719
$ws.Range("L21").Value = 0.2117089459261621
$ws.Range("M21").Value = 0.2751501752986059
$ws.Range("B22").Value = 1.199227172909929
$ws.Range("C22").Value = 0.05217264008339839
$ws.Range("D22").Value = 0.4185586659220064
$ws.Range("E22").Value = 0.0882051058079476
$ws.Range("F22").Value = 6.195492625782776
$ws.Range("K22").Value = 0.7148312888898829
$ws.Range("L22").Value = 0.2143448110729622
$ws.Range("M22").Value = 0.2791633996231155
$ws.Range("B23").Value = 1.189051103595034
$ws.Range("C23").Value = 0.05120197458670361
$ws.Range("D23").Value = 0.4126294595663467
$ws.Range("E23").Value = 0.08785669317097344
$ws.Range("F23").Value = 6.090513180697087
$ws.Range("K23").Value = 0.705749120307388
$ws.Range("L23").Value = 0.2129298932229631
$ws.Range("M23").Value = 0.2770056444772706
$ws.Range("B24").Value = 1.152256039308185
$ws.Range("C24").Value = 0.04748664263112801
$ws.Range("D24").Value = 0.3900455902670785
$ws.Range("E24").Value = 0.0865575685737543
$ws.Range("F24").Value = 5.693597949439209
$ws.Range("K24").Value = 0.672666277722243
$ws.Range("L24").Value = 0.2077490700347084
$ws.Range("M24").Value = 0.2691805821507174
$ws.Range("B25").Value = 1.116275059608654
$ws.Range("C25").Value = 0.04338883934056526
$ws.Range("D25").Value = 0.3653962687009908
$ws.Range("E25").Value = 0.08519964973832117
$ws.Range("F25").Value = 5.266750071844967
$ws.Range("K25").Value = 0.6397619673392967
$ws.Range("L25").Value = 0.2025399877624849
$ws.Range("M25").Value = 0.2614784560245944